$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from AC1 onto the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record (Wins/Losses/Ties) for every player row
$lastRow = 39
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 83   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 79   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}

$wb.Save()
